$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Offers" (G2) and "Discount" (H2) values from the data row.
# This also drops the now-unused "12"/"10" shared strings automatically.
$ws.Range("G2:H2").Clear()

# Move the active selection to K2 (slider/scrolling focus cell).
$ws.Range("K2").Select()
